# Append 7 new blood-pressure measurement rows (rows 122-128) below the
# existing data in "Arkusz1", matching the author's exported rows, then
# update the sheet's scrolled position / selection to reflect the new
# bottom of the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Columns: Date (A), Time (B), Systolic (C), Diastolic (D), Pulse (E)
$newRows = @(
    @(45973, 0.54166666666666663, 109, 65, 84),
    @(45973, 0.66666666666666663, 112, 63, 74),
    @(45973, 0.79166666666666663, 116, 75, 70),
    @(45973, 0.91666666666666663, 114, 69, 78),
    @(45974, 0.41666666666666669, 101, 63, 73),
    @(45974, 0.54166666666666663, 112, 60, 80),
    @(45974, 0.66666666666666663, 106, 65, 77)
)

$lastExistingRow = 121
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Carry the formatting (number formats / borders / alignment) of the
    # last existing data row down onto the freshly appended one first, so
    # new cells pick up the same style indices as A121:E121.
    $ws.Range("A$lastExistingRow`:E$lastExistingRow").Copy() | Out-Null
    $ws.Range("A$r`:E$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
}

$excel.CutCopyMode = $false

$lastRow = $startRow + $newRows.Count - 1   # 128

# Scroll the sheet so the new rows are visible, and leave the selection on
# the cell below the freshly entered data, like Excel does after typing.
$excel.ActiveWindow.ScrollRow = $lastExistingRow
$ws.Range("A$lastRow").Select() | Out-Null
